$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full target table for the doctyp_code master-data sheet (header + 35 data rows).
$rows = @(
    ,@("doctyp_code", "doccat_code", "lang_code", "is_active", "cr_by", "cr_dtimes")
    ,@("CIN", "POI", "ara", $true, "superadmin", "now()")
    ,@("RNC", "POA", "ara", $true, "superadmin", "now()")
    ,@("COR", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC001", "POI", "ara", $true, "superadmin", "now()")
    ,@("CRN", "POR", "ara", $true, "superadmin", "now()")
    ,@("COB", "POB", "ara", $false, "superadmin", "now()")
    ,@("DOC001", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC002", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC003", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC004", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC005", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC006", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC007", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC008", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC009", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC010", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC011", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC012", "POI", "ara", $true, "superadmin", "now()")
    ,@("DOC001", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC013", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC014", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC015", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC004", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC005", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC006", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC016", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC017", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC018", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC008", "POA", "ara", $true, "superadmin", "now()")
    ,@("DOC024", "POR", "ara", $true, "superadmin", "now()")
    ,@("DOC025", "POR", "ara", $true, "superadmin", "now()")
    ,@("DOC026", "POR", "ara", $true, "superadmin", "now()")
    ,@("DOC001", "POR", "ara", $true, "superadmin", "now()")
    ,@("DOC027", "POR", "ara", $true, "superadmin", "now()")
    ,@("DOC028", "POR", "ara", $true, "superadmin", "now()")
)

$rowCount = $rows.Count
$colCount = 6
$data = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$targetRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, $colCount))
$targetRange.Value = $data

# Selection lands past the last populated column, matching a manual paste
# (mirrors the original file's "whole unused columns" selection style).
$ws.Range("G1:XFD1048576").Select()
